# Restore C10 on the active sheet from 18 to 1 (value-only edit; the
# surrounding style/formatting attributes for the cell are left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
